# Generate Report for Handback
# Adds a new handback row (file 0f366ddc-864c-400c-8b92-96137923b8a2.md) to the
# Overview / zh-cn / de-de sheets, mirroring the existing 647b1ad8... row.

$wb = $excel.ActiveWorkbook

$missing = [System.Reflection.Missing]::Value

# ---------------------------------------------------------------------------
# Sheet "Overview"  (table3 -> A1:G2 => A1:G3)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = "0f366ddc-864c-400c-8b92-96137923b8a2.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d9dec251e76b39edc38cf024ba04dd9b32024a01/e2e/0f366ddc-864c-400c-8b92-96137923b8a2.md", $missing, $missing, "e2e\0f366ddc-864c-400c-8b92-96137923b8a2.md") | Out-Null
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Range("G3").Value = "2016-10-18 03:53:46"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"  (table1 -> A1:P2 => A1:P3)
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/d9dec251e76b39edc38cf024ba04dd9b32024a01/e2e/0f366ddc-864c-400c-8b92-96137923b8a2.md", $missing, $missing, "0f366ddc-864c-400c-8b92-96137923b8a2.md") | Out-Null
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "True"
$wsZhCn.Range("G3").Value = "0f366ddc-864c-400c-8b92-96137923b8a2.d9dec251e76b39edc38cf024ba04dd9b32024a01.zh-cn.xlf"
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("H3").Value = "2016-10-18 03:53:30"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/d9dec251e76b39edc38cf024ba04dd9b32024a01/e2e/0f366ddc-864c-400c-8b92-96137923b8a2.md", $missing, $missing, "0f366ddc-864c-400c-8b92-96137923b8a2.md") | Out-Null
$wsZhCn.Range("J3").Value = "0f366ddc-864c-400c-8b92-96137923b8a2.d9dec251e76b39edc38cf024ba04dd9b32024a01.zh-cn.xlf"
$wsZhCn.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("K3").Value = "2016-10-18 03:54:26"
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("O3").Value = "False"

# ---------------------------------------------------------------------------
# Sheet "de-de"  (table2 -> A1:P2 => A1:P3)
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/d9dec251e76b39edc38cf024ba04dd9b32024a01/e2e/0f366ddc-864c-400c-8b92-96137923b8a2.md", $missing, $missing, "0f366ddc-864c-400c-8b92-96137923b8a2.md") | Out-Null
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "True"
$wsDeDe.Range("G3").Value = "0f366ddc-864c-400c-8b92-96137923b8a2.d9dec251e76b39edc38cf024ba04dd9b32024a01.de-de.xlf"
$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("H3").Value = "2016-10-18 03:53:46"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/d9dec251e76b39edc38cf024ba04dd9b32024a01/e2e/0f366ddc-864c-400c-8b92-96137923b8a2.md", $missing, $missing, "0f366ddc-864c-400c-8b92-96137923b8a2.md") | Out-Null
$wsDeDe.Range("J3").Value = "0f366ddc-864c-400c-8b92-96137923b8a2.d9dec251e76b39edc38cf024ba04dd9b32024a01.de-de.xlf"
$wsDeDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("K3").Value = "2016-10-18 03:54:49"
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("O3").Value = "False"

Write-Host "Handback row added to Overview, zh-cn and de-de sheets."
